# Added Power Button Control Circuit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 previously held "TPS61232DRC " (A14) and "MUST HAVE FIXED 5V out" (B14).
# That part (power design note) is replaced by a new "Power Control Transistor"
# entry placed a couple of rows below, mirroring the TPS61230DRCR hyperlink entry
# above it. Clear out the old A14/B14 content, keep A14's wrapped-text style.
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()

# The row shrinks back down to the single-line height now that it is empty.
$ws.Rows("14").RowHeight = 13.8

# New entries for the Power Control Transistor (row 16) and its datasheet /
# purchase link (row 17), leaving row 15 blank as a visual spacer just like
# the other gaps in this sheet (rows 4, 6, 8, 12, 15).
$ws.Range("A16").Value2 = "Power Control Tranistor"
$ws.Rows("16").RowHeight = 15

$ws.Range("A17").Value2 = "https://www.digikey.com/en/products/detail/goford-semiconductor/G6N02L/13664832"
$ws.Rows("17").RowHeight = 13.8

# Update the saved selection to match the author's cursor position after the edit.
$ws.Range("C23").Select()
